$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "B2" "Bitcoin"
Set-TextCell $ws "C2" "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell $ws "D2" "22.623.76"
Set-TextCell $ws "E2" "  -4.11%  "

Set-TextCell $ws "B3" "Ethereum"
Set-TextCell $ws "C3" "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell $ws "D3" "1.584.85"
Set-TextCell $ws "E3" "  -4.15%  "

Set-TextCell $ws "B4" "TetherUSD"
Set-TextCell $ws "C4" "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell $ws "D4" "1.014"
Set-TextCell $ws "E4" "  +1.19%  "

Set-TextCell $ws "B5" "USDC"
Set-TextCell $ws "C5" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell $ws "D5" "1.013"
Set-TextCell $ws "E5" "  +1.07%  "

Set-TextCell $ws "B6" "BNB"
Set-TextCell $ws "C6" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell $ws "D6" "294.91"
Set-TextCell $ws "E6" "  -2.40%  "

Set-TextCell $ws "B7" "XRP"
Set-TextCell $ws "C7" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell $ws "D7" "0.3651"
Set-TextCell $ws "E7" "  -3.80%  "

Set-TextCell $ws "B8" "OKB"
Set-TextCell $ws "C8" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D8" "49.55"
Set-TextCell $ws "E8" "  -2.35%  "

Set-TextCell $ws "B9" "Cardano"
Set-TextCell $ws "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell $ws "D9" "0.3326"
Set-TextCell $ws "E9" "  -6.75%  "

Set-TextCell $ws "B10" "Polygon"
Set-TextCell $ws "C10" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell $ws "D10" "1.167"
Set-TextCell $ws "E10" "  -5.31%  "

Set-TextCell $ws "B11" "Dogecoin"
Set-TextCell $ws "C11" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell $ws "D11" "0.07549"
Set-TextCell $ws "E11" "  -7.15%  "

Set-TextCell $ws "B12" "BinanceUSD"
Set-TextCell $ws "C12" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell $ws "D12" "1.013"
Set-TextCell $ws "E12" "  +1.18%  "

Set-TextCell $ws "B13" "Solana"
Set-TextCell $ws "C13" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell $ws "D13" "20.99"
Set-TextCell $ws "E13" "  -5.53%  "

Set-TextCell $ws "B14" "Polkadot"
Set-TextCell $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D14" "5.988"
Set-TextCell $ws "E14" "  -7.06%  "

Set-TextCell $ws "B15" "Chainlink"
Set-TextCell $ws "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D15" "6.863"
Set-TextCell $ws "E15" "  -8.08%  "

Set-TextCell $ws "B16" "WrappedEther"
Set-TextCell $ws "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D16" "1.588.83"
Set-TextCell $ws "E16" "  -3.50%  "

Set-TextCell $ws "B17" "ShibaInu"
Set-TextCell $ws "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D17" "0.00001133"
Set-TextCell $ws "E17" "  -6.16%  "

Set-TextCell $ws "B18" "Litecoin"
Set-TextCell $ws "C18" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws "D18" "90.17"
Set-TextCell $ws "E18" "  -7.42%  "

Set-TextCell $ws "B19" "TRON"
Set-TextCell $ws "C19" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D19" "0.06830"
Set-TextCell $ws "E19" "  -2.38%  "

Set-TextCell $ws "B20" "Dai"
Set-TextCell $ws "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D20" "1.014"
Set-TextCell $ws "E20" "  +1.20%  "

Set-TextCell $ws "B21" "Uniswap"
Set-TextCell $ws "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws "D21" "6.219"
Set-TextCell $ws "E21" "  -8.15%  "

Set-TextCell $ws "B22" "Avalanche"
Set-TextCell $ws "C22" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws "D22" "16.29"
Set-TextCell $ws "E22" "  -7.35%  "

Set-TextCell $ws "B23" "Cosmos"
Set-TextCell $ws "C23" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D23" "11.86"
Set-TextCell $ws "E23" "  -6.68%  "

Set-TextCell $ws "B24" "WrappedBTC"
Set-TextCell $ws "C24" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws "D24" "22.634.87"
Set-TextCell $ws "E24" "  -4.11%  "

Set-TextCell $ws "B25" "Toncoin"
Set-TextCell $ws "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D25" "2.385"
Set-TextCell $ws "E25" "  -4.31%  "

Set-TextCell $ws "B26" "LidoDAOToken"
Set-TextCell $ws "C26" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D26" "2.886"
Set-TextCell $ws "E26" "  -2.18%  "

Set-TextCell $ws "B27" "EthereumClassic"
Set-TextCell $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D27" "19.81"
Set-TextCell $ws "E27" "  -6.32%  "

Set-TextCell $ws "B28" "Monero"
Set-TextCell $ws "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D28" "144.66"
Set-TextCell $ws "E28" "  -4.99%  "

Set-TextCell $ws "B29" "HuobiToken"
Set-TextCell $ws "C29" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws "D29" "5.009"
Set-TextCell $ws "E29" "  -3.83%  "

Set-TextCell $ws "B30" "BitcoinCash"
Set-TextCell $ws "C30" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D30" "125.70"
Set-TextCell $ws "E30" "  -5.89%  "

Set-TextCell $ws "B31" "WrappedliquidstakedEther2.0"
Set-TextCell $ws "C31" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell $ws "D31" "1.769.31"
Set-TextCell $ws "E31" "  -3.13%  "

Set-TextCell $ws "B32" "WEMIXTOKEN"
Set-TextCell $ws "C32" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D32" "2.084"
Set-TextCell $ws "E32" "  -3.35%  "

Set-TextCell $ws "B33" "Filecoin"
Set-TextCell $ws "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D33" "6.295"
Set-TextCell $ws "E33" "  -11.14%  "

Set-TextCell $ws "B34" "FraxShare"
Set-TextCell $ws "C34" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D34" "10.42"
Set-TextCell $ws "E34" "  -12.26%  "

Set-TextCell $ws "B35" "ImmutableX"
Set-TextCell $ws "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D35" "0.9169"
Set-TextCell $ws "E35" "  -11.62%  "

Set-TextCell $ws "B36" "Stellar"
Set-TextCell $ws "C36" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D36" "0.08508"
Set-TextCell $ws "E36" "  -2.82%  "

Set-TextCell $ws "B37" "VeChain"
Set-TextCell $ws "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D37" "0.02524"
Set-TextCell $ws "E37" "  -8.25%  "

Set-TextCell $ws "B38" "Algorand"
Set-TextCell $ws "C38" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws "D38" "0.2275"
Set-TextCell $ws "E38" "  -7.74%  "

Set-TextCell $ws "B39" "InternetComputer(DFINITY)"
Set-TextCell $ws "C39" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D39" "5.427"
Set-TextCell $ws "E39" "  -10.13%  "

Set-TextCell $ws "B40" "Hedera"
Set-TextCell $ws "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D40" "0.06354"
Set-TextCell $ws "E40" "  -8.07%  "

Set-TextCell $ws "B41" "TrustWalletToken"
Set-TextCell $ws "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D41" "1.270"
Set-TextCell $ws "E41" "  -4.30%  "

Set-TextCell $ws "B42" "Aptos"
Set-TextCell $ws "C42" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws "D42" "11.85"
Set-TextCell $ws "E42" "  -10.21%  "

Set-TextCell $ws "B43" "TheSandbox"
Set-TextCell $ws "C43" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws "D43" "0.6343"
Set-TextCell $ws "E43" "  -9.02%  "

Set-TextCell $ws "B44" "EnergySwap"
Set-TextCell $ws "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D44" "14.47"
Set-TextCell $ws "E44" "  -9.02%  "

Set-TextCell $ws "B45" "Frax"
Set-TextCell $ws "C45" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws "D45" "1.012"
Set-TextCell $ws "E45" "  +1.05%  "

Set-TextCell $ws "B46" "PancakeSwap"
Set-TextCell $ws "C46" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D46" "3.845"
Set-TextCell $ws "E46" "  -2.44%  "

Set-TextCell $ws "B47" "Decentraland"
Set-TextCell $ws "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws "D47" "0.5880"
Set-TextCell $ws "E47" "  -9.31%  "

Set-TextCell $ws "B48" "NEARProtocol"
Set-TextCell $ws "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D48" "2.096"
Set-TextCell $ws "E48" "  -8.10%  "

Set-TextCell $ws "B49" "Quant"
Set-TextCell $ws "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws "D49" "122.43"
Set-TextCell $ws "E49" "  -3.71%  "

Set-TextCell $ws "B50" "Cronos"
Set-TextCell $ws "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D50" "0.07296"
Set-TextCell $ws "E50" "  -7.72%  "

Set-TextCell $ws "B51" "EOS"
Set-TextCell $ws "C51" "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextCell $ws "D51" "1.169"
Set-TextCell $ws "E51" "  -1.75%  "
